# Updated calculations for 57MHz pixel clock.
#
# B1 ("clock (mhz)") used to hold the formula =113.75/2 (56.875 MHz).
# It is now a hard-coded 57 MHz literal. Every other populated cell on
# the sheet is a formula that (directly or transitively) depends on B1,
# so writing the new literal and letting Excel recalculate reproduces
# the entire cascade of updated values shown in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace the formula in B1 with the literal value 57 (drops the old
# =113.75/2 formula entirely).
$ws.Range("B1").Value = 57

# The saved view's selection moved from J26 to B2.
$ws.Range("B2").Select()
